# top_remaining_CA.xlsx refresh: new scrape pass (2019-03-07 -> 2019-03-12).
# Most rows just get a bumped "TOP PRIZES REMAINING" count and a refreshed
# "LAST SCRAPE DATE". A handful of price tiers also had their GAME NAME /
# GAME NUMBER rows reordered between scrapes (the underlying scraper re-sorts
# by remaining prize count), which is why some C/D values below swap between
# two adjacent rows rather than just changing in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- GAME NAME / GAME NUMBER / TOP PRIZES REMAINING updates ---
$ws.Range("E2").Value = 29
$ws.Range("E3").Value = 30

# $2.00 Games tier reshuffle: AMAZING ACES drops behind MONOPOLY / IT TAKES 2
$ws.Range("C8").Value = "MONOPOLY"
$ws.Range("D8").Value = 1348
$ws.Range("E8").Value = 25
$ws.Range("C9").Value = "IT TAKES 2"
$ws.Range("D9").Value = 1340
$ws.Range("E9").Value = 28
$ws.Range("C10").Value = "AMAZING ACES"
$ws.Range("D10").Value = 1330
$ws.Range("E10").Value = 9

$ws.Range("E17").Value = 36
$ws.Range("E21").Value = 9
$ws.Range("E22").Value = 9

# $5.00 Games tier reshuffle: MEGA CROSSWORD drops behind 20X THE CASH / CROSSWORD CONNECT
$ws.Range("C23").Value = "20X THE CASH"
$ws.Range("D23").Value = 1315
$ws.Range("C24").Value = "CROSSWORD CONNECT"
$ws.Range("D24").Value = 1341
$ws.Range("E24").Value = 11
$ws.Range("C25").Value = "MEGA CROSSWORD"
$ws.Range("D25").Value = 1308
$ws.Range("E25").Value = 2

# ROAD TO RICHES / MONOPOLY swap places
$ws.Range("C26").Value = "MONOPOLY"
$ws.Range("D26").Value = 1349
$ws.Range("E26").Value = 19
$ws.Range("C27").Value = "ROAD TO RICHES"
$ws.Range("D27").Value = 1311
$ws.Range("E27").Value = 5

# SILVER RICHES / $59M PRIZE POOL swap places
$ws.Range("C29").Value = "$59M PRIZE POOL"
$ws.Range("D29").Value = 1328
$ws.Range("E29").Value = 4
$ws.Range("C30").Value = "SILVER RICHES"
$ws.Range("D30").Value = 1336
$ws.Range("E30").Value = 5

$ws.Range("E36").Value = 5
$ws.Range("E37").Value = 10
$ws.Range("E38").Value = 45744

# $20.00 Games tier: MONOPOLY / $5,000,000 Spectacular swap places
$ws.Range("C44").Value = "$5,000,000 Spectacular"
$ws.Range("D44").Value = 1304
$ws.Range("E44").Value = 1
$ws.Range("C45").Value = "MONOPOLY"
$ws.Range("D45").Value = 1351
$ws.Range("E45").Value = 4

# TRIPLE JACKPOT / JACKPOT FORTUNE swap places
$ws.Range("C46").Value = "JACKPOT FORTUNE"
$ws.Range("D46").Value = 1338
$ws.Range("C47").Value = "TRIPLE JACKPOT"
$ws.Range("D47").Value = 1300

# --- LAST SCRAPE DATE updates ---
# Force text formatting on the whole date column first so Excel doesn't
# auto-convert the "yyyy-mm-dd" strings into date serial numbers, then
# clear the formatting again afterwards so no extra cell styles are left
# behind (matches the source file, where this column is plain text).
$dateRange = $ws.Range("F2:F50")
$dateRange.NumberFormat = "@"

# Rows that were last scraped 2019-03-07 move to the new 2019-03-12 pass.
# (Rows 7, 12, 19 and 20 were already on an older 2019-02-19 pass and are
# untouched; row 47 is the exception below, following JACKPOT FORTUNE's
# old 2019-03-07 stamp down into TRIPLE JACKPOT's new slot.)
$bumpedRows = @(2,3,4,5,6,8,9,10,11,13,14,15,16,17,18,21,22,23,24,25,26,27,28,
                29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,48,49,50)
foreach ($r in $bumpedRows) {
    $ws.Range("F$r").Value = "2019-03-12"
}
$ws.Range("F47").Value = "2019-02-19"

$dateRange.ClearFormats()
